$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '71.185.36'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +6.16%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.666.21'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +5.91%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '596.10'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '195.64'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +4.70%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.651'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +2.67%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.659.92'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +5.91%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.182'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +5.51%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.677'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +4.98%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '59.08'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +4.94%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000295'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +6.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.02'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +6.78%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.255.79'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +5.81%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.02'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +7.00%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.674.07'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +5.78%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '71.236.25'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.85'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +5.84%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.122'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +3.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.07'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +5.81%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '494.45'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.89%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '19.05'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +13.88%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.50%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.54'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '91.99'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +2.57%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.18'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +8.40%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.58'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +6.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.69'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +7.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '33.18'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +5.88%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.95'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +11.69%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +9.77%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '633.71'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +5.73%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '12.34'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +5.55%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '66.02'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +3.17%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '40.87'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +12.28%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0₃0842'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +11.81%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.415'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +8.49%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.71%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.62'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +2.51%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.332.10'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.17'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +9.18%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.87'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +14.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0456'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +6.63%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.92'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +5.15%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +3.38%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.31'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.25%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.33'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +7.40%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.11%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.12%  '
